$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column E ("reviews_count"), shifting columns F:K left to E:J
$ws.Columns.Item(5).Delete()
